$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 10

# Delete rows 3-5 entirely (Sofia Vicente, Miguel Gutierrez, Perrito)
$ws.Range("A3:C5").EntireRow.Delete()
